$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timing results for rows 2-13 (columns A-D); column E (dataset
# labels) is unchanged.
$values = @{
    2  = @(422.99361062049871, 1167.1050772666929, 106.89730668067931, 527.38269352912903)
    3  = @(26.143944025039669, 452.81845545768738, 14.689358711242679, 55.088619947433472)
    4  = @(337.89333248138428, 870.56675434112549, 99.246004104614258, 637.76248741149902)
    5  = @(80.430171728134155, 480.26403498649597, 48.461355924606323, 113.8541820049286)
    6  = @(48.260196924209588, 384.91414165496832, 15.31342434883118, 1685.406959533691)
    7  = @(116.8032686710358, 697.29553961753845, 56.120332479476929, 99.909819364547729)
    8  = @(2001.4262361526489, 900.11027455329895, 66.310784816741943, 1708.473667621613)
    9  = @(444.96363496780401, 822.1752917766571, 81.83507251739502, 7861.6777367591858)
    10 = @(113.13900852203371, 488.70281982421881, 33.191932916641242, 3112.2439227104192)
    11 = @(258.34844613075262, 583.76041150093079, 118.5075278282166, 1219.428186655045)
    12 = @(2323.9030058383942, 2421.3472683429718, 1169.3939740657811, 13507.267176151279)
    13 = @(672.82793235778809, 825.22512316703796, 648.71337890625, 7525.5196301937103)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 1).Value = $rowValues[0]
    $ws.Cells.Item($row, 2).Value = $rowValues[1]
    $ws.Cells.Item($row, 3).Value = $rowValues[2]
    $ws.Cells.Item($row, 4).Value = $rowValues[3]
}

# Row 13 loses its bordered/bold "header-like" style that had been applied
# (s="2" -> default), while keeping its values/text (column E untouched).
$ws.Range("A13:E13").ClearFormats()
